$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.357.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.89%  "

$ws.Range("D3").Value = "'3.609.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.86%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'240.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.74%  "

$ws.Range("D6").Value = "'638.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.08%  "

$ws.Range("E7").Value = "  +10.07%  "

$ws.Range("E8").Value = "  +6.86%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("E10").Value = "  +9.71%  "

$ws.Range("D11").Value = "'3.607.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.88%  "

$ws.Range("D12").Value = "'43.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.38%  "

$ws.Range("E13").Value = "  +5.29%  "

$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.63%  "

$ws.Range("D15").Value = "'4.276.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.98%  "

$ws.Range("D16").Value = "'96.258.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.96%  "

$ws.Range("E17").Value = "  +6.08%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'8.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.50%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.612.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.95%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'13.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +24.75%  "

$ws.Range("D21").Value = "'18.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.01%  "

$ws.Range("D22").Value = "'0.500"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.07%  "

$ws.Range("D23").Value = "'516.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.10%  "

$ws.Range("E24").Value = "  +2.44%  "

$ws.Range("E25").Value = "  +12.63%  "

$ws.Range("E26").Value = "  +11.27%  "

$ws.Range("D27").Value = "'93.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.49%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'3.798.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.87%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'12.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'3.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +20.23%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'11.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.55%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.143"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.17%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "'0.182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.43%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'0.995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.09%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'30.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.00%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.569"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.67%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'575.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.83%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.21%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'1.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.55%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.929"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.91%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "'1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.39%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0431"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.63%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'23.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'5.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.40%  "

$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").Value = "'3.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.82%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'53.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.16%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.76%  "
